# data preparing tranformers creation
# Remove the stray "shape"/"type" columns (F, G) that were added during
# anomaly testing, fix up a few Rho (D) readings, and clear the now-unused
# helper values in column E.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the extraneous "shape" / "type" columns (F1:G2 held the only data in
# those columns: headers "shape"/"type" and values "C"/"CB2P").
$ws.Range("F1:G2").ClearContents()

# Column E only ever held two stray scalar flags (E2, E10) - clear them but
# keep their cell styling in place.
$ws.Range("E2").ClearContents()
$ws.Range("E10").ClearContents()

# Correct a handful of Rho (column D) readings.
$ws.Range("D2").Value = 230
$ws.Range("D5").Value = 125
$ws.Range("D9").Value = 168

# Move the active selection to reflect where editing left off.
$ws.Range("D16").Select() | Out-Null

# The chart anchored over columns H:P shifted slightly (same vertical
# extent, a touch narrower) after the data cleanup above - nudge it to
# match.
$co = $ws.ChartObjects().Item(1)
$co.Top = $co.Top + 11.25
$co.Left = $co.Left + 19.4375
$co.Width = $co.Width - 10.4375
